$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Add default formats for dates and date times..." note that
# lived at G12:I12 - its text is being replaced/relocated into a new task row.
$ws.Range("G12:I12").Clear()

# Row 21 ("Specify baseStyleName and inline styles...") switches from the
# "Neutral" (yellow) style to the "Good" (green) style, now that the feature
# has been implemented.
$ws.Range("B21:E21").Style = "Good"

# Insert 5 new blank rows above the old row 23, pushing the existing task
# rows (23-27) down to (28-32) so four new "Issue" tasks can be added above
# them.
$ws.Range("23:27").Insert()

# New task rows 23-26 (only columns B:E are styled/used for these rows).
$ws.Range("B23").Value = "Add quick pivot tests (see trello)."
$ws.Range("B23:E23").Style = "Neutral"

$ws.Range("B24").Value = "Issue 9:  Support changing PivotDataGroup caption"
$ws.Range("B24:E24").Style = "Neutral"

$ws.Range("B25").Value = "Issue 3:  defineCalculation(filters) should be able to replace filters, not just combine with AND"
$ws.Range("B25:E25").Style = "Neutral"

$ws.Range("B26").Value = "Issue 1:  calculationType ""value"" should work with totals"
$ws.Range("B26:E26").Style = "Neutral"

# New row 33 - the relocated "Add note to the docs..." task, styled like the
# row above it (row 32, the bold/neutral "In NEWS.md..." row).
$ws.Range("B33").Value = "Add note to the docs that for excel export, outputting values as rawValue for dates/posixct will cause a number to appear."
$ws.Range("B33:F33").Style = "Neutral"

# Restore the selection to match the recorded cursor position.
[void]$ws.Range("G24:G25").Select()
